# excel_writer: also include totals for the balance columns
#
# The "Gesamtergebnis" (totals) sheet previously left the balance columns
# (Startguthaben/C and Endsaldo/D) in the "Total" row as "N/A" text. Now
# they should be included in the totals as numeric 0, just like every
# other column on that row.
#
# This also updates the saved selection/active-sheet state to match what
# a user would see after making this change in the workbook (the totals
# sheet becomes the active tab, with C3:D3 selected there; the other two
# sheets keep their cursor on A2).

$wb = $excel.ActiveWorkbook

$wsDaily   = $wb.Worksheets.Item("Tagesergebnisse")
$wsMonthly = $wb.Worksheets.Item("Monatsergebnisse")
$wsTotal   = $wb.Worksheets.Item("Gesamtergebnis")

# --- Core data fix: balance column totals on the "Total" row (row 3) ---
# Replace the "N/A" placeholders in the balance columns with real numeric
# totals (0, matching every other totals column on that row).
$wsTotal.Range("C3").Value = 0
$wsTotal.Range("D3").Value = 0

# The row's rendered height shrinks slightly now that it holds numbers
# instead of text.
$wsTotal.Rows.Item(3).RowHeight = 13.8

# --- Selection / active sheet bookkeeping ---
# Leave the cursor on A2 on the other two sheets...
$wsDaily.Range("A2").Select()
$wsMonthly.Range("A2").Select()

# ...and make the totals sheet active, with the newly-filled-in balance
# cells selected.
$wsTotal.Activate()
$wsTotal.Range("C3:D3").Select()
